# Add descriptions/titles for the DMI Identifiant Local Fabricant extension.
$wb = $excel.ActiveWorkbook

# --- "Metadata" sheet: Property/Value summary table ---
$meta = $wb.Worksheets.Item("Metadata")

# Title
$meta.Range("B5").Value = "DMI Identifiant Local Fabricant"

# Date (bumped to reflect the edit)
$meta.Range("B8").Value = "2026-02-25T08:15:31+00:00"

# Description
$meta.Range("B12").Value = "Extension créée dans ce volet pour représenter l'identifiant local fabricant."

# --- "Elements" sheet: StructureDefinition element table, root Extension row (row 2) ---
$elements = $wb.Worksheets.Item("Elements")

# Short
$elements.Range("L2").Value = "DMI Identifiant Local Fabricant"

# Definition
$elements.Range("M2").Value = "Extension créée dans ce volet pour représenter l'identifiant local fabricant."

# Mapping: RIM Mapping -> cleared
$elements.Range("AK2").Value = ""
